$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DM_Stat (C) and P_Value (D) figures after correction to Diebold Mariano test
$ws.Range("C2").Value = -1.382148465108401
$ws.Range("D2").Value = 0.1759395920405216

$ws.Range("C3").Value = 0.4647498477817668
$ws.Range("D3").Value = 0.645073258614882

$ws.Range("C4").Value = 0.7980702673087269
$ws.Range("D4").Value = 0.4303688737379554

$ws.Range("C5").Value = 1.627499632067555
$ws.Range("D5").Value = 0.1128643614560494

$ws.Range("C6").Value = 1.2174976765984
$ws.Range("D6").Value = 0.2317976813596569

$ws.Range("C7").Value = 1.428292060967515
$ws.Range("D7").Value = 0.1623301767175664

$ws.Range("C8").Value = 2.030197382656884
$ws.Range("D8").Value = 0.0502184890067241

$ws.Range("C9").Value = 0.08232161111900922
$ws.Range("D9").Value = 0.9348737076812657

$ws.Range("C10").Value = 1.025508311082208
$ws.Range("D10").Value = 0.3123645780081572

$ws.Range("C11").Value = 1.34699599118755
$ws.Range("D11").Value = 0.1868914601856302
